# NIT-9012426536.xlsx : "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# - Add a new payment-period row (period 2202) at the bottom of the data table,
#   re-order the existing period rows so the most recent period is listed first,
#   add the new period 2507, update the worker-count/period-count summary, and
#   bump the overdue-balance total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new data row right after the current last data row (56) so the
#    table grows from 41 periods (rows 16-56) to 42 periods (rows 16-57).
$ws.Rows(57).Insert()

# 2) Fix up the row formatting: row 57 (the new bottom row) should carry the
#    "last row" (thicker bottom border) formatting that row 56 used to have,
#    and row 56 should fall back to the regular interior-row formatting (like
#    row 55), since it is no longer the last row of the table.
$ws.Range("B56:J56").Copy()
$ws.Range("B57:J57").PasteSpecial(-4122)
$ws.Range("B55:J55").Copy()
$ws.Range("B56:J56").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Fill in the new bottom row's data (same worker/doc/salary as every other
#    row, oldest period 2202).
$ws.Range("B57").Value = "PE"
$ws.Range("C57").Value = "776402523111963"
$ws.Range("D57").Value = "JOSE LUIS PEREZ"
$ws.Range("E57").Value = "2202"
$ws.Range("F57").Value = 40000
$ws.Range("G57").Value = 1000000

# 4) Re-order the period column so it now runs most-recent-first (2507 down
#    to 2202) across all 42 data rows (16-57).
$ws.Range("E16").Value = "2507"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2504"
$ws.Range("E20").Value = "2503"
$ws.Range("E21").Value = "2502"
$ws.Range("E22").Value = "2501"
$ws.Range("E23").Value = "2412"
$ws.Range("E24").Value = "2411"
$ws.Range("E25").Value = "2410"
$ws.Range("E26").Value = "2409"
$ws.Range("E27").Value = "2408"
$ws.Range("E28").Value = "2407"
$ws.Range("E29").Value = "2406"
$ws.Range("E30").Value = "2405"
$ws.Range("E31").Value = "2404"
$ws.Range("E32").Value = "2403"
$ws.Range("E33").Value = "2402"
$ws.Range("E34").Value = "2401"
$ws.Range("E35").Value = "2312"
$ws.Range("E36").Value = "2311"
$ws.Range("E37").Value = "2310"
$ws.Range("E38").Value = "2309"
$ws.Range("E39").Value = "2308"
$ws.Range("E40").Value = "2307"
$ws.Range("E41").Value = "2306"
$ws.Range("E42").Value = "2305"
$ws.Range("E43").Value = "2304"
$ws.Range("E44").Value = "2303"
$ws.Range("E45").Value = "2302"
$ws.Range("E46").Value = "2301"
$ws.Range("E47").Value = "2212"
$ws.Range("E48").Value = "2211"
$ws.Range("E49").Value = "2210"
$ws.Range("E50").Value = "2209"
$ws.Range("E51").Value = "2208"
$ws.Range("E52").Value = "2207"
$ws.Range("E53").Value = "2206"
$ws.Range("E54").Value = "2205"
$ws.Range("E55").Value = "2204"
$ws.Range("E56").Value = "2203"
$ws.Range("E57").Value = "2202"

# 5) Update the summary cells: overdue balance total and number of periods.
$ws.Range("E11").Value = 1680000
$ws.Range("F13").Value = 42
